$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A82").Value = "2023-12-07 19:13:41"
$ws.Range("B82").Value = 0.004000000000000001

$ws.Range("A83").Value = "2023-12-07 19:14:06"
$ws.Range("B83").Value = 0.002

$ws.Range("A84").Value = "2023-12-07 19:14:17"
$ws.Range("B84").Value = 0.0008
